$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# remain Text (matching the source inlineStr cells) instead of being
# auto-converted to a Number by Excel's type inference. We set NumberFormat
# to Text ("@") before assigning, then reset the cell style back to "Normal"
# so no residual style index is left on the cell (it keeps storing as text).

$ws.Range("D2").Value = "26.906.02"
$ws.Range("E2").Value = "  +1.14%  "

$ws.Range("D3").Value = "1.842.97"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4751"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.84%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3674"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9256"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.916.79"
$ws.Range("E12").Value = "  +5.01%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07651"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.310"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.400"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("D20").Value = "26.927.77"
$ws.Range("E20").Value = "  +1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.049"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.923"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.14"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.002"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.945"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08854"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.289"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7483"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.171"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.754"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.480"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.090"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05260"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01948"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.961"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5209"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.965"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1513"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4726"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.58%  "

$ws.Range("E48").Value = "  +3.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06028"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8855"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.35%  "
